$wb = $excel.ActiveWorkbook

# --- Update Portugal sheet's selection (it is no longer the active tab) ---
$portugal = $wb.Worksheets.Item("Portugal")
[void]$portugal.Activate()
[void]$portugal.Range("D17").Select()

# --- Create the new "Slovakia" sheet by copying "Portugal" (keeps styles/merges/col widths) ---
[void]$portugal.Copy($null, $portugal)
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Remove the rows for "P32AR" and "P32DR" (not present in the Slovakia market)
[void]$slovakia.Rows.Item(16).Delete()
[void]$slovakia.Rows.Item(16).Delete()

# Reset the (now unnecessary) tall row height inherited from Portugal's rows 3-5
[void]$slovakia.Range("A3:A5").EntireRow.AutoFit()

# Fill in the market-specific values.
# NOTE: B4 is written before B2 so that the new shared-string entries are
# appended in the same order as the source workbook (NGC id, then market name).
$slovakia.Range("B4").Value = "NGC-2930/T3222"
$slovakia.Range("B2").Value = "Slovakia Market"

# Extend the blank "B" column formatting down through the data rows, matching Portugal's layout
[void]$slovakia.Range("B3").Copy()
[void]$slovakia.Range("B6:B19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Place the selection/active-cell and make this the visible tab, as in the final workbook
[void]$slovakia.Range("B4").Select()
